$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "git commit calificacion hasta p23": a new grading/scoring field "l1" is
# inserted right before the existing "nota_iniciativa" field. Column Q used
# to hold the "nota_iniciativa" header; it now becomes "l1", and a brand new
# column R is appended holding the "nota_iniciativa" header (with the same
# 0-flag data pattern the old column had).

$lastRow = $ws.UsedRange.Rows.Count

# 1. Rename header of column Q from "nota_iniciativa" to "l1"
$ws.Range("Q1").Value = "l1"

# 2. Add the new column R header "nota_iniciativa", matching the other
#    header cells' formatting (bold, bordered).
$ws.Range("R1").Value = "nota_iniciativa"
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Fill R2:R<lastRow> with 0, mirroring columns P and Q's existing data.
$ws.Range("R2:R$lastRow").Value = 0
